$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Maria Claudia%Alzamora%NULL%0,                           Tania%Paredes%NULL%2,                           Tania%Paredes%NULL%0,                           David%Caceres%NULL%1,                           Camille M.%Webb%NULL%1,                           Luis M.%Valdez%NULL%1,                           Mauricio%La Rosa%NULL%1]'
$ws.Range("E3").Value = '[Francesco%Catellani%NULL%1,                           Andrea%Coscione%NULL%1,                           Riccardo%D’Ambrosi%NULL%1,                           Luca%Usai%NULL%1,                           Claudio%Roscitano%NULL%1,                           Gennaro%Fiorentino%NULL%1]'
$ws.Range("E4").Value = '[Rong%Chen%NULL%0,                           Yuan%Zhang%NULL%1,                           Lei%Huang%NULL%0,                           Bi-heng%Cheng%NULL%1,                           Zhong-yuan%Xia%NULL%1,                           Qing-tao%Meng%mengqingtao2018@126.com%1]'
$ws.Range("E5").Value = '[Lan%Dong%xref no email%3,       Jinhua%Tian%xref no email%3,       Songming%He%xref no email%3,       Chuchao%Zhu%xref no email%3,       Jian%Wang%xref no email%0,       Chen%Liu%xref no email%3,       Jing%Yang%xref no email%3]'
$ws.Range("E6").Value = '[Yin%Du%NULL%1,                           Long%Wang%NULL%1,                           Gang%Wu%NULL%1,                           Xiaoming%Lei%NULL%1,                           Wei%Li%NULL%0,                           Jianrui%Lv%lvjianrui@163.com%1]'
$ws.Range("E7").Value = '[E%Ferrazzi%NULL%1,                           L%Frigerio%NULL%1,                           V%Savasi%NULL%1,                           P%Vergani%NULL%1,                           F%Prefumo%NULL%1,                           S%Barresi%NULL%2,                           S%Barresi%NULL%0,                           S%Bianchi%NULL%1,                           E%Ciriello%NULL%1,                           F%Facchinetti%NULL%1,                           MT%Gervasi%NULL%1,                           E%Iurlaro%NULL%1,                           A%Kustermann%NULL%1,                           G%Mangili%NULL%1,                           F%Mosca%NULL%1,                           L%Patanè%NULL%1,                           D%Spazzini%NULL%1,                           A%Spinillo%NULL%1,                           G%Trojano%NULL%1,                           M%Vignali%NULL%1,                           A%Villa%NULL%1,                           GV%Zuccotti%NULL%1,                           F%Parazzini%Fabio.parazzini@unimi.it%1,                           I%Cetin%NULL%1]'
$ws.Range("E8").Value = '[Michael S.%Firstenberg%msfirst@gmail.com%1,                           Matthew%Libby%mnlibby@gmail.com%1,                           Michael%Ochs%Michael.Ochs@usap.com%1,                           Jennifer%Hanna%hannaj21@hotmail.com%1,                           Julie E.%Mangino%Jemangino1@gmail.com%1,                           Joseph%Forrester%Joe.Forrester@ccpcmd.com%1]'
$ws.Range("E9").Value = '[Yunhe%Gao%NULL%1,                           Hongqing%Xi%NULL%1,                           Lin%Chen%NULL%0]'
$ws.Range("E10").Value = '[Sebastian%Gidl\u00f6f%xref no email%0,       Julia%Savchenko%xref no email%0,       Thomas%Brune%xref no email%0,       Hampus%Josefsson%xref no email%0]'
$ws.Range("E11").Value = '[Hao%He%NULL%1,                           Shuai%Zhao%NULL%2,                           Linlin%Han%NULL%1,                           Qi%Wang%NULL%1,                           Haifa%Xia%NULL%1,                           Xin%Huang%NULL%1,                           Shanglong%Yao%NULL%2,                           Jiapeng%Huang%NULL%2,                           Xiangdong%Chen%xiangdongchen2013@163.com%3]'
$ws.Range("E12").Value = '[Dong Hwan%Lee%NULL%0,                           Jihyang%Lee%NULL%2,                           Jihyang%Lee%NULL%0,                           Eunju%Kim%NULL%2,                           Eunju%Kim%NULL%0,                           Kyeongyoon%Woo%NULL%2,                           Kyeongyoon%Woo%NULL%0,                           Hak Youle%Park%NULL%2,                           Hak Youle%Park%NULL%0,                           Jihyun%An%NULL%2,                           Jihyun%An%NULL%0]'
$ws.Range("E13").Value = '[Yang%Li%NULL%0,                           Ruihong%Zhao%NULL%2,                           Shufa%Zheng%NULL%1,                           Xu%Chen%NULL%0,                           Jinxi%Wang%NULL%1,                           Xiaoli%Sheng%NULL%1,                           Jianying%Zhou%NULL%1,                           Hongliu%Cai%NULL%1,                           Qiang%Fang%NULL%0,                           Fei%Yu%NULL%1,                           Jian%Fan%NULL%1,                           Kaijin%Xu%NULL%1,                           Yu%Chen%NULL%2,                           Jifang%Sheng%NULL%0]'
$ws.Range("E14").Value = '[Dawei%Lu%NULL%0,                           Lin%Sang%NULL%1,                           Shihua%Du%NULL%1,                           Tao%Li%NULL%0,                           Yange%Chang%NULL%1,                           Xiu‐An%Yang%yangxiuan07@mails.ucas.edu.cn%1]'
$ws.Range("E15").Value = '[Joana%Lyra%xref no email%1,       Rita%Valente%xref no email%1,       Marta%Ros\u00e1rio%xref no email%1,       Mariana%Guimar\u00e3es%xref no email%1]'
$ws.Range("E16").Value = '[Bobin%Mi%xref no email%1,       Lang%Chen%xref no email%1,       Adriana C.%Panayi%xref no email%0,       Jing%Liu%xref no email%0,       Faqi%Cao%xref no email%1,       Chenchen%Yan%xref no email%1,       Hang%Xue%xref no email%0,       Yiqiang%Hu%xref no email%1,       Liangcong%Hu%xref no email%1,       Ze%Lin%xref no email%0,       Xudong%Xie%xref no email%0,       Ranyang%Tao%xref no email%1,       Yun%Sun%xref no email%1,       Wu%Zhou%xref no email%0,       Yuan%Xiong%xref no email%1,       Guohui%Liu%xref no email%0]'
$ws.Range("E17").Value = '[James E%Archer%NULL%1,                           Abdulrahman%Odeh%NULL%1,                           Simon%Ereidge%NULL%1,                           Hosni Khairy%Salem%NULL%1,                           Gareth P%Jones%NULL%1,                           Anne%Gardner%NULL%1,                           Shiva S%Tripathi%NULL%1,                           Anne%Gregg%NULL%1,                           Reuben%Jeganathan%NULL%1,                           Muhammad Harris%Siddique%NULL%1,                           Celestino P%Lombardi%NULL%1,                           Benjamin%Martin%NULL%1,                           Cosimo Alex%Leo%NULL%1,                           Debashis%Dass%NULL%1,                           Gregorio%Di Franco%NULL%1,                           Long R%Jiao%NULL%1,                           Giulio Maria%Mari%NULL%1,                           Luis-Cristobal%Capitan-Morales%NULL%1,                           Tara M%Connelly%NULL%1,                           Ammar%Alanbuki%NULL%1,                           Armando%De Virgilio%NULL%1,                           Clare%Schilling%NULL%1,                           Carlos%San Miguel Mendez%NULL%1,                           Gauri%Kulkarni%NULL%1,                           Kulsoom%Nizami%NULL%1,                           Stewart%Walsh%NULL%1,                           Harry%Dean%NULL%1,                           Miguel%Ruiz-Marin%NULL%1,                           Rory%Houston%NULL%1,                           Mario%Trompetto%NULL%1,                           David%Chrastek%NULL%1,                           Vasileios%Kouritas%NULL%1,                           Maria%Cannoletta%NULL%1,                           Francesco%Rosato%NULL%1,                           Manish V%Kaushal%NULL%1,                           Paulo Matos%Costa%NULL%1,                           Hannah H%Elkadi%NULL%1,                           Jack R%Johnstone%NULL%1,                           Esmee%Irvine%NULL%1,                           Maria R%Alvarez%NULL%1,                           Carlo%Corbellini%NULL%1,                           Gowtham S%Venkatesan%NULL%1,                           Olga%Mateo-Sierra%NULL%1,                           Carolina%Martínez-Pérez%NULL%1,                           Javier%Serrano González%NULL%1,                           Miguel Õngel%Hernández Bartolomé%NULL%1,                           David%Díaz Pérez%NULL%1,                           María%Gutiérrez Samaniego%NULL%1,                           Pablo%Galindo Jara%NULL%1,                           Neil%Sharma%NULL%1]'
$ws.Range("E18").Value = '[Song%Limin%xref no email%1,       Xiao%Weimin%xref no email%1,       Ling%Ken%xref no email%1,       Yao%Shanglong%xref no email%1,       Chen%Xiangdong%xref no email%1]'
$ws.Range("E19").Value = '[Mingyang%Sun%NULL%1,                           Guoting%Xu%NULL%1,                           Yong%Yang%NULL%1,                           Yuan%Tao%NULL%1,                           May%Pian-Smith%NULL%1,                           Vandana%Madhavan%NULL%1,                           Zhongcong%Xie%zxie@mgh.harvard.edu%1,                           Jiaqiang%Zhang%jqzhang@henu.edu.cn%1]'
$ws.Range("E20").Value = '[Xiaotong%Wang%NULL%0,                           Zhiqiang%Zhou%NULL%2,                           Jianping%Zhang%NULL%2,                           Fengfeng%Zhu%NULL%2,                           Yongyan%Tang%NULL%2,                           Xinghua%Shen%luckydrx@163.com%2]'
$ws.Range("E21").Value = '[Haifa%Xia%NULL%1,                           Shuai%Zhao%NULL%0,                           Zhouyang%Wu%NULL%1,                           Huilin%Luo%NULL%1,                           Cheng%Zhou%NULL%1,                           Xiangdong%Chen%NULL%0]'
$ws.Range("E22").Value = '[Hui%Zeng%xref no email%1,       Chen%Xu%xref no email%1,       Junli%Fan%xref no email%1,       Yueting%Tang%xref no email%1,       Qiaoling%Deng%xref no email%2,       Wei%Zhang%xref no email%1,       Xinghua%Long%xref no email%1]'
$ws.Range("E23").Value = '[Zhi-Jiang%Zhang%NULL%1,                           Xue-Jie%Yu%NULL%1,                           Tao%Fu%NULL%1,                           Yu%Liu%NULL%2,                           Yu%Liu%NULL%0,                           Yan%Jiang%NULL%1,                           Bing Xiang%Yang%NULL%1,                           Yongyi%Bi%NULL%2,                           Yongyi%Bi%NULL%0]'
$ws.Range("E24").Value = '[Shuai%Zhao%NULL%1,                           Ken%Ling%NULL%1,                           Hong%Yan%NULL%1,                           Liang%Zhong%NULL%1,                           Xiaohong%Peng%NULL%1,                           Shanglong%Yao%NULL%0,                           Jiapeng%Huang%NULL%0,                           Xiangdong%Chen%NULL%0]'
$ws.Range("E25").Value = '[Qi%Zhong%NULL%1,                           Yin Y.%Liu%NULL%1,                           Qiong%Luo%NULL%1,                           Yu F.%Zou%NULL%1,                           Hai X.%Jiang%NULL%1,                           Hui%Li%NULL%0,                           Jing J.%Zhang%NULL%1,                           Zhen%Li%NULL%0,                           Xin%Yang%NULL%1,                           Min%Ma%NULL%1,                           Li J.%Tang%NULL%1,                           Ying Y.%Chen%chenyingy@whu.edu.cn%1,                           Feng%Zheng%fengzheng@whu.edu.cn%1,                           Jian J.%Ke%1219628972@qq.com%1,                           Zong Z.%Zhang%zhangzz@whu.edu.cn%1]'
$ws.Range("E27").Value = '[Benvenuti%M. R.%coreGivesNoEmail%1,                          Bonardelli%S.%coreGivesNoEmail%1,                          Calza%S.%coreGivesNoEmail%1,                          Casiraghi%A.%coreGivesNoEmail%1,                          Castelli%F.%coreGivesNoEmail%1,                          Doglietto%F. (ORCID:0000-0002-7438-0734)%coreGivesNoEmail%1,                          Domenicucci%M.%coreGivesNoEmail%1,                          Fontanella%M. M.%coreGivesNoEmail%1,                          Gheza%F.%coreGivesNoEmail%1,                          Lussardi%G. L.%coreGivesNoEmail%1,                          Maroldi%R.%coreGivesNoEmail%1,                          Milano%G.%coreGivesNoEmail%1,                          Panciani%P. P.%coreGivesNoEmail%1,                          Portolani%N.%coreGivesNoEmail%1,                          Rasulo%F. A.%coreGivesNoEmail%1,                          Saraceno%G.%coreGivesNoEmail%1,                          Signorini%L.%coreGivesNoEmail%1,                          Vecchiarelli%L.%coreGivesNoEmail%1,                          Vezzoli%M.%coreGivesNoEmail%1,                          Zanin%L.%coreGivesNoEmail%1]'
$ws.Range("C26").Value = 'Management of Spine Trauma in COVID-19 Pandemic: A Preliminary Report.'
$ws.Range("D26").Value = 'Background COVID-19 was first identified in Iran in February 2020 and since then it spread rapidly through all over the country and soon after that it was reported as a pandemic.
 The current study presents a preliminary report of spine trauma management during COVID-19 pandemic.
 Methods A cross sectional study was designed to evaluate patients admitted for vertebral fractures with diagnosis of COVID-19 infection on February and March 2020. Analysis was made based on clinical and laboratory data along with the imaging findings from chest HRCT.
 Results Seven patients with spine trauma including five males and two females ranging from 14 to 59 years were diagnosed for COVID-19 infection through CT-scan findings.
 Except one, all other patients were asymptomatic for COVID-19 at the time of admission.
 In three cases the COVID diagnosis was made the day after arrival and in others after 10, 14 and 35 days.
 Five patients were treated surgically among whom four were admitted to ICU soon after the surgery.
 The mean ICU stay for operated patients were eight days and the mean hospital stay was 22.6 days.
 Conclusion Proper diagnosis of COVID-19 is the keystone to protect both patients and health care providers.
 During the pandemic all admitted patients should be screened for COVID-19 infection.
 Unnecessary procedures for spine trauma patients should be avoided in order to reduce complications related to surgery and to preserve ICU beds.
'
$ws.Range("E26").Value = '[ Mohammadreza%Chehrassan%elasticNoEmail%1,   Adel% Ebrahimpour%elasticNoEmail%1,   Hasan% Ghandhari%elasticNoEmail%1,   Morteza% Sanei Taheri%elasticNoEmail%1,   Bahador% Athari%elasticNoEmail%1,   Mehrdad% Sadighi%elasticNoEmail%1,   Meisam% Jafari KafiAbadi%elasticNoEmail%1,   Amin% Karami%elasticNoEmail%1,   Alireza% Zali%elasticNoEmail%1]'
$ws.Range("F26").Value = '10.22038/abjs.2020.47882.2368'
$ws.Range("G26").Value = 'ELASTIC'
$ws.Range("H26").Value = '2020-04-01'
